$wb = $excel.ActiveWorkbook

# --- View_Print sheet: replace formula in A1 with text, set B1 text, move selection ---
$wsPrint = $wb.Worksheets.Item("View_Print")
$wsPrint.Activate()
$wsPrint.Range("A1").Value = "Log 1"
$wsPrint.Range("B1").Value = "Log 2"
$wsPrint.Range("A1:O1").Select()
$excel.ActiveCell = $wsPrint.Range("O1")

# --- Log sheet: fix header row (restore "Log N" names), drop column P ---
$wsLog = $wb.Worksheets.Item("Log")
$wsLog.Activate()
$wsLog.Range("A10").Value = "Log 1"
$wsLog.Range("B10").Value = "Log 2"
$wsLog.Range("C10").Value = "Log 3"
$wsLog.Range("D10").Value = "Log 4"
$wsLog.Range("E10").Value = "Log 5"
$wsLog.Range("F10").Value = "Log 6"
$wsLog.Range("G10").Value = "Log 7"
$wsLog.Range("H10").Value = "Log 8"
$wsLog.Range("I10").Value = "Log 9 "
$wsLog.Range("J10").Value = "Log 10"
$wsLog.Range("K10").Value = "Log 11"
$wsLog.Range("L10").Value = "Log 12"
$wsLog.Range("M10").Value = "Log 13"
$wsLog.Range("N10").Value = "Log 14"
$wsLog.Range("O10").Value = "Log 15"
$wsLog.Range("P10").Clear()
$wsLog.Range("B13").Select()

# --- Input sheet: clear A2 value, move selection to A2, keep as active tab ---
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Activate()
$wsInput.Range("A2").ClearContents()
$wsInput.Range("A2").Select()

$wb.Save()
